# Adds units to headers across the input sheets (Added units to input files)
$wb = $excel.ActiveWorkbook

# --- Fleet sheet: p_max_ch -> p_max_ch (kW), p_max_ds -> p_max_ds (kW) ---
# This must happen FIRST so the new shared strings are appended in the same
# order as the target workbook.
$wsFleet = $wb.Worksheets.Item("Fleet")
$wsFleet.Range("C1").Value = "p_max_ch (kW)"
$wsFleet.Range("D1").Value = "p_max_ds (kW)"
[void]$wsFleet.Range("A1").Select()

# --- Cluster1 / Cluster2 / Cluster3: cu_p_ch_max -> cu_p_ch_max (kW), cu_p_ds_max -> cu_p_ds_max (kW) ---
$wsCluster1 = $wb.Worksheets.Item("Cluster1")
$wsCluster1.Range("B1").Value = "cu_p_ch_max (kW)"
$wsCluster1.Range("C1").Value = "cu_p_ds_max (kW)"
[void]$wsCluster1.Range("B1:C1").Select()

$wsCluster2 = $wb.Worksheets.Item("Cluster2")
$wsCluster2.Range("B1").Value = "cu_p_ch_max (kW)"
$wsCluster2.Range("C1").Value = "cu_p_ds_max (kW)"
[void]$wsCluster2.Range("B1:C1").Select()

$wsCluster3 = $wb.Worksheets.Item("Cluster3")
$wsCluster3.Range("B1").Value = "cu_p_ch_max (kW)"
$wsCluster3.Range("C1").Value = "cu_p_ds_max (kW)"
[void]$wsCluster3.Range("B1:C1").Select()

# --- Capacity1 / Capacity2 / Capacity3: LB -> LB (kW), UB -> UB (kW) ---
$wsCapacity1 = $wb.Worksheets.Item("Capacity1")
$wsCapacity1.Range("B1").Value = "LB (kW)"
$wsCapacity1.Range("C1").Value = "UB (kW)"
[void]$wsCapacity1.Range("B1:C1").Select()

$wsCapacity2 = $wb.Worksheets.Item("Capacity2")
$wsCapacity2.Range("B1").Value = "LB (kW)"
$wsCapacity2.Range("C1").Value = "UB (kW)"
[void]$wsCapacity2.Range("B1:C1").Select()

$wsCapacity3 = $wb.Worksheets.Item("Capacity3")
$wsCapacity3.Range("B1").Value = "LB (kW)"
$wsCapacity3.Range("C1").Value = "UB (kW)"
[void]$wsCapacity3.Range("B1:C1").Select()

# --- Price sheet: Price -> Price (per/kWh) ---
$wsPrice = $wb.Worksheets.Item("Price")
$wsPrice.Range("B1").Value = "Price (per/kWh)"
[void]$wsPrice.Range("B1").Select()

# Re-select the Price sheet last (it is the tab that remains active/selected).
[void]$wsPrice.Activate()
